# Add Russia, Finland and Hungary market test-data sheets, mirroring the
# existing per-country "Repeaters" sheets.
#
# Russia & Hungary use the same 19-row layout as "Austria" (no MZXSDR240
# row); Finland uses the 20-row layout as "Denmark" (includes MZXSDR240).

$wb = $excel.ActiveWorkbook

$austriaSrc = $wb.Worksheets.Item("Austria")
$denmarkSrc = $wb.Worksheets.Item("Denmark")

# --- Russia (copied from Austria, placed after Denmark) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$austriaSrc.Copy($null, $lastSheet)
$russia = $wb.Worksheets.Item($wb.Worksheets.Count)
$russia.Name = "Russia"
$russia.Range("B4").Value = "NGC-2929/T2910"
$russia.Range("B2").Value = "Russia Market"

# --- Finland (copied from Denmark, placed after Russia) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmarkSrc.Copy($null, $lastSheet)
$finland = $wb.Worksheets.Item($wb.Worksheets.Count)
$finland.Name = "Finland"
$finland.Range("B4").Value = "NGC-3130/T2887"
$finland.Range("B2").Value = "Finland Market"

# --- Hungary (copied from Austria, placed after Finland) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$austriaSrc.Copy($null, $lastSheet)
$hungary = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungary.Name = "Hungary"
$hungary.Range("B4").Value = "NGC-3104/T2979"
$hungary.Range("B2").Value = "Hungary Market"

# Narrow column B on the new sheets to fit the short "NGC-xxxx/Txxxx" /
# "<Country> Market" values instead of the long user-story text (matches
# the auto-fitted width Excel computes once the long user-story text is
# gone from that column).
foreach ($sheet in @($russia, $finland, $hungary)) {
    $sheet.Columns.Item(2).ColumnWidth = 14.29
    $sheet.Rows.Item(2).RowHeight = 28.8
    $sheet.Rows.Item(3).RowHeight = 28.8
    $sheet.Rows.Item(4).RowHeight = 28.8
    $sheet.Rows.Item(5).RowHeight = 28.8
}

# Match the selection state left behind on each new sheet.
$russia.Range("A1:D19").Select() | Out-Null
$finland.Range("A1:D19").Select() | Out-Null

# Hungary is the new active/selected sheet.
$hungary.Select()
$hungary.Range("I12").Select() | Out-Null

Write-Output "done"
